$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: ypddjuio@yomail.info -> anhhuy9900@yopmail.com
$hl2 = $ws.Range("C2").Hyperlinks.Item(1)
$hl2.Address = "mailto:anhhuy9900@yopmail.com"
$hl2.TextToDisplay = "anhhuy9900@yopmail.com"

# C3: kawnlyiw@supere.ml -> anhhuy9901@yopmail.com
$hl3 = $ws.Range("C3").Hyperlinks.Item(1)
$hl3.Address = "mailto:anhhuy9901@yopmail.com"
$hl3.TextToDisplay = "anhhuy9901@yopmail.com"

# C4: pythonprocourse2@gmail.com -> nhahuy29051990@gmail.com
$hl4 = $ws.Range("C4").Hyperlinks.Item(1)
$hl4.Address = "mailto:nhahuy29051990@gmail.com"
$hl4.TextToDisplay = "nhahuy29051990@gmail.com"
